# Add a new day ("07-nov") of data to "Dataframe ST.xlsx"
#
# 1. Sheet3 holds the raw per-product values for the newest business day in
#    column B (rows 20:36); Sheet3!C2:C18 (and Sheet1!CB:CC, which reference
#    Sheet3 via VLOOKUP) recalculate automatically from those raw values.
# 2. Sheet1 gets one more date column appended (CJ) holding the same VLOOKUP
#    result as CB/CC, mirroring the existing CD..CI "frozen" daily columns.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1) Update the raw per-product figures on Sheet3 (rows 20-36, column B)
#    for the new day. Everything that depends on these (Sheet3!C and
#    Sheet1!CB/CC) recalculates automatically.
# ---------------------------------------------------------------------
$newRawValues = @{
    20 = 11.965584732690733
    21 = 0
    22 = 0
    23 = 5.085721638017886
    24 = 6.1995403354425003
    25 = 0
    26 = 5.630951299561973
    27 = 15.523260272918142
    28 = 4.0474094049327913
    29 = 3.8146778192564836
    30 = 8.9280179955007064
    31 = 15.422350855629702
    32 = 4.2310265131564799
    33 = 12.82072511631856
    34 = 7.4757604770686132
    35 = 5.2086136085173402
    36 = 48.916269942172214
}

foreach ($r in $newRawValues.Keys) {
    $ws3.Cells.Item($r, 2).Value = $newRawValues[$r]
}

$excel.CalculateFull()

# ---------------------------------------------------------------------
# 2) Append the new "07-nov" date column (CJ) to Sheet1.
# ---------------------------------------------------------------------
$ws1.Range("CJ1").Value = "07-nov"
$ws1.Range("CJ1").NumberFormat = "@"
$ws1.Range("CJ2:CJ18").NumberFormat = "0"

for ($r = 2; $r -le 18; $r++) {
    $ws1.Cells.Item($r, 88).Value = $ws1.Cells.Item($r, 80).Value()
}

# ---------------------------------------------------------------------
# 3) Leave the selection where the author left it after the edit.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("CI22").Select()
